$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G3').Value = 'Tumakuru (Tumkur)'
$ws.Range('G4').Value = 'Ballari (Bellary)'
$ws.Range('G5').Value = 'Ballari (Bellary)'
$ws.Range('G6').Value = 'Tumakuru (Tumkur)'
$ws.Range('G7').Value = 'Ballari (Bellary)'
$ws.Range('G8').Value = 'Tumakuru (Tumkur)'
$ws.Range('G9').Value = 'Ballari (Bellary)'
$ws.Range('G10').Value = 'Tumakuru (Tumkur)'
$ws.Range('G11').Value = 'Tumakuru (Tumkur)'
$ws.Range('G12').Value = 'Ballari (Bellary)'
$ws.Range('G14').Value = 'Chitradurga'
$ws.Range('G15').Value = 'Ballari (Bellary)'
$ws.Range('G16').Value = 'Tumakuru (Tumkur)'
$ws.Range('G17').Value = 'Tumakuru (Tumkur)'
$ws.Range('G18').Value = 'Tumakuru (Tumkur)'
$ws.Range('G20').Value = 'Tumakuru (Tumkur)'
$ws.Range('G21').Value = 'Tumakuru (Tumkur)'
$ws.Range('G22').Value = 'Tumakuru (Tumkur)'
$ws.Range('G23').Value = 'Tumakuru (Tumkur)'
$ws.Range('G24').Value = 'Tumakuru (Tumkur)'
$ws.Range('G25').Value = 'Ballari (Bellary)'
$ws.Range('G26').Value = 'Ballari (Bellary)'
$ws.Range('G27').Value = 'Udupi'
$ws.Range('G28').Value = 'Tumakuru (Tumkur)'
$ws.Range('G29').Value = 'Tumakuru (Tumkur)'
$ws.Range('G30').Value = 'Ballari (Bellary)'
$ws.Range('G31').Value = 'Tumakuru (Tumkur)'
$ws.Range('G32').Value = 'Tumakuru (Tumkur)'
$ws.Range('G33').Value = 'Tumakuru (Tumkur)'
$ws.Range('G34').Value = 'Ballari (Bellary)'
$ws.Range('G35').Value = 'Tumakuru (Tumkur)'
$ws.Range('G36').Value = 'Ballari (Bellary)'
$ws.Range('G37').Value = 'Tumakuru (Tumkur)'
$ws.Range('G38').Value = 'Tumakuru (Tumkur)'
$ws.Range('G39').Value = 'Chitradurga'
$ws.Range('G40').Value = 'Chitradurga'
$ws.Range('G41').Value = 'Tumakuru (Tumkur)'
$ws.Range('G42').Value = 'Ballari (Bellary)'
$ws.Range('G43').Value = 'Ballari (Bellary)'
$ws.Range('G44').Value = 'Tumakuru (Tumkur)'
$ws.Range('G45').Value = 'Tumakuru (Tumkur)'
$ws.Range('G46').Value = 'Tumakuru (Tumkur)'
$ws.Range('G47').Value = 'Chitradurga'
$ws.Range('G48').Value = 'Tumakuru (Tumkur)'
$ws.Range('G49').Value = 'Tumakuru (Tumkur)'
$ws.Range('G50').Value = 'Ballari (Bellary)'
$ws.Range('G51').Value = 'Ballari (Bellary)'
$ws.Range('G52').Value = 'Ballari (Bellary)'
$ws.Range('G53').Value = 'Tumakuru (Tumkur)'
$ws.Range('G54').Value = 'Ballari (Bellary)'
$ws.Range('G55').Value = 'Tumakuru (Tumkur)'
$ws.Range('G56').Value = 'Tumakuru (Tumkur)'
$ws.Range('G57').Value = 'Chitradurga'
$ws.Range('G58').Value = 'Chitradurga'
$ws.Range('G59').Value = 'Tumakuru (Tumkur)'
$ws.Range('G60').Value = 'Ballari (Bellary)'
$ws.Range('G61').Value = 'Tumakuru (Tumkur)'
$ws.Range('G62').Value = 'Tumakuru (Tumkur)'
$ws.Range('G64').Value = 'Tumakuru (Tumkur)'
$ws.Range('G65').Value = 'Tumakuru (Tumkur)'
$ws.Range('G66').Value = 'Ballari (Bellary)'
$ws.Range('G67').Value = 'Tumakuru (Tumkur)'
$ws.Range('G68').Value = 'Tumakuru (Tumkur)'
$ws.Range('G69').Value = 'Tumakuru (Tumkur)'
$ws.Range('G70').Value = 'Tumakuru (Tumkur)'
$ws.Range('G71').Value = 'Udupi'
$ws.Range('G72').Value = 'Tumakuru (Tumkur)'
$ws.Range('G73').Value = 'Tumakuru (Tumkur)'
$ws.Range('G74').Value = 'Tumakuru (Tumkur)'
